$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157674908638
$ws.Range("B1").Value = 2.255504369735718
$ws.Range("C1").Value = 3.033737182617188
$ws.Range("D1").Value = 1.429672241210938
$ws.Range("E1").Value = 1.026046752929688
